$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BOM corrections (swap wrong resistor/capacitor values for the correct parts) ---
# Row 4: R3,R4 -> now 4.3k 100mW resistors (ERJ3GEYJ432V)
$ws.Range("E4").Value = "ERJ3GEYJ432V"

# Row 3: R1,R2 -> now 2.2k 100mW resistors (ERJ3GEYJ222V)
$ws.Range("E3").Value = "ERJ3GEYJ222V"

$ws.Range("B4").Value = "4.3k 100mW"
$ws.Range("B3").Value = "2.2k 100mW"

# Row 5: was mislabeled as R2/1.5uF -> now correctly C1,C2 1uF 25V capacitors (CL10A105KA8NNNC)
$ws.Range("E5").Value = "CL10A105KA8NNNC"
$ws.Range("B5").Value = "1uF 25V"
$ws.Range("D5").Value = "C1, C2"

# --- Stray formatted (empty) cells left next to the edited rows, Verdana 8pt FF333333 ---
[void]$ws.Range("G3:G4").Select()
$ws.Range("G3").Font.Size = 8
$ws.Range("G3").Font.Color = 3355443
$ws.Range("G3").Font.Name = "Verdana"

$ws.Range("G4").Font.Size = 8
$ws.Range("G4").Font.Color = 3355443
$ws.Range("G4").Font.Name = "Verdana"

# --- Restore the selection left by the author after editing ---
[void]$ws.Range("B3:F5").Select()
